$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B23").Value = 6330
$ws.Range("C23").Value = 1002
$ws.Range("D23").Value = 5922319
$ws.Range("E23").Value = 935.5954186413902
$ws.Range("F23").Value = 8.613589567604674
$ws.Range("G23").Value = 4.266389177939645
$ws.Range("H23").Value = 26.90013719990066
